$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511688082662454"
$ws1.Range("B2").Value = "go_stims-16511688082242475.csv"
$ws1.Range("B3").Value = "GNG_stims-1651168808249284.csv"
$ws1.Range("B4").Value = "go_stims-16511688082502444.csv"
$ws1.Range("B5").Value = "GNG_stims-16511688082652793.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-165116881049084"
$ws2.Range("B2").Value = "OB-16511688093191643.csv"
$ws2.Range("B3").Value = "ZB-match_7-16511688088817482.csv"
$ws2.Range("B4").Value = "TB-16511688104752204.csv"
$ws2.Range("B5").Value = "ZB-match_4-16511688090970514.csv"
$ws2.Range("B6").Value = "OB-16511688091763167.csv"
$ws2.Range("B7").Value = "TB-16511688098492103.csv"
$ws2.Range("B8").Value = "OB-16511688092382138.csv"
$ws2.Range("B9").Value = "ZB-match_3-16511688083462481.csv"
$ws2.Range("B10").Value = "TB-16511688097867117.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-165116881049084"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651168810537726"
$ws4.Range("B2").Value = "MM_stims-1651168810506471.csv"
$ws4.Range("B3").Value = "ZM_stims-165116881049084.csv"
$ws4.Range("B4").Value = "MM_stims-16511688105220623.csv"
$ws4.Range("B5").Value = "ZM_stims-1651168810506471.csv"
$ws4.Range("B6").Value = "MM_stims-1651168810537726.csv"
$ws4.Range("B7").Value = "ZM_stims-16511688105220623.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511688106002223"
$ws5.Range("B2").Value = "SAT_stims-16511688105533123.csv"
$ws5.Range("B3").Value = "SAT_stims-1651168810537726.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511688105845978.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511688105689871.csv"
